$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 204.11111
$ws.Range("I19").Value = 277.6
$ws.Range("J19").Value = 112.25
$ws.Range("K19").Value = 277.6
$ws.Range("L19").Value = 112.25
$ws.Range("M19").Value = -102.6
$ws.Range("N19").Value = -462.25

$ws.Range("H29").Value = 596.4286
$ws.Range("J29").Value = 791.6667
$ws.Range("L29").Value = 2375.0001
$ws.Range("N29").Value = -2937.0001

$ws.Range("H58").Value = 672.7778
$ws.Range("I58").Value = 40
$ws.Range("J58").Value = 1179
$ws.Range("K58").Value = 120
$ws.Range("L58").Value = 3537
$ws.Range("M58").Value = 30
$ws.Range("N58").Value = -3837

$ws.Range("H70").Value = 3432.1428
$ws.Range("I70").Value = 3075
$ws.Range("K70").Value = 9225
$ws.Range("M70").Value = -8955

$ws.Range("H73").Value = 3432.1428
$ws.Range("I73").Value = 3075
$ws.Range("K73").Value = 9225
$ws.Range("M73").Value = -8289

$ws.Range("H80").Value = 212.55556
$ws.Range("I80").Value = 251.85715
$ws.Range("J80").Value = 75
$ws.Range("K80").Value = 755.5714499999999
$ws.Range("L80").Value = 225
$ws.Range("M80").Value = 242.4285500000001
$ws.Range("N80").Value = -2221

$ws.Range("H83").Value = 212.55556
$ws.Range("I83").Value = 251.85715
$ws.Range("J83").Value = 75
$ws.Range("K83").Value = 2266.71435
$ws.Range("L83").Value = 675
$ws.Range("M83").Value = 2725.28565
$ws.Range("N83").Value = -10659

$ws.Range("H113").Value = 1907.5
$ws.Range("I113").Value = 2422.5
$ws.Range("K113").Value = 2422.5
$ws.Range("M113").Value = 831.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 526
$ws.Range("I2").Value = 504.33334
$ws.Range("K2").Value = 504.33334
$ws.Range("M2").Value = -391.33334

$ws.Range("H15").Value = 4000
$ws.Range("I15").Value = 4000
$ws.Range("K15").Value = 4000
$ws.Range("M15").Value = -3650

$ws.Range("H25").Value = 4784.1
$ws.Range("I25").Value = 1787.5
$ws.Range("J25").Value = 6781.8335
$ws.Range("K25").Value = 1787.5
$ws.Range("L25").Value = 6781.8335
$ws.Range("M25").Value = -1385.5
$ws.Range("N25").Value = -7585.8335

$ws.Range("H32").Value = 5519.625
$ws.Range("I32").Value = 4480.3335
$ws.Range("K32").Value = 4480.3335
$ws.Range("M32").Value = -4193.3335

$ws.Range("H35").Value = 2019.5
$ws.Range("I35").Value = 2019.5
$ws.Range("K35").Value = 2019.5
$ws.Range("M35").Value = -1613.5

$ws.Range("H116").Value = 526
$ws.Range("I116").Value = 504.33334
$ws.Range("K116").Value = 504.33334
$ws.Range("M116").Value = 1789.66666

$ws.Range("H132").Value = 3300
$ws.Range("I132").Value = 3300
$ws.Range("K132").Value = 9900
$ws.Range("M132").Value = -7370

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 526
$ws.Range("I3").Value = 504.33334
$ws.Range("K3").Value = 504.33334
$ws.Range("M3").Value = -390.33334

$ws.Range("H25").Value = 1044.5
$ws.Range("I25").Value = 1044.5
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1044.5
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = -809.5

$ws.Range("H26").Value = 48333.332
$ws.Range("I26").Value = 15000
$ws.Range("J26").Value = 65000
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 65000
$ws.Range("M26").Value = -14708
$ws.Range("N26").Value = -65584

$ws.Range("H28").Value = 59999
$ws.Range("J28").Value = 59999
$ws.Range("L28").Value = 59999
$ws.Range("N28").Value = -60587

$ws.Range("H36").Value = 10690.25
$ws.Range("I36").Value = 9431.857
$ws.Range("K36").Value = 9431.857
$ws.Range("M36").Value = -8897.857

$ws.Range("H37").Value = 3231.3333
$ws.Range("I37").Value = 894
$ws.Range("K37").Value = 894
$ws.Range("M37").Value = -757

$ws.Range("H96").Value = 2085.6
$ws.Range("I96").Value = 2085.6
$ws.Range("K96").Value = 2085.6
$ws.Range("M96").Value = 660.4000000000001

$ws.Range("H99").Value = 1766.8182
$ws.Range("I99").Value = 1693.5
$ws.Range("K99").Value = 1693.5
$ws.Range("M99").Value = -195.5

$ws.Range("H105").Value = 1733.8
$ws.Range("I105").Value = 1539.5
$ws.Range("K105").Value = 1539.5
$ws.Range("M105").Value = 207.5

$ws.Range("H107").Value = 50778.25
$ws.Range("I107").Value = 67289.664
$ws.Range("K107").Value = 67289.664
$ws.Range("M107").Value = -65369.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 21978.334
$ws.Range("J28").Value = 21978.334
$ws.Range("L28").Value = 21978.334
$ws.Range("N28").Value = -22468.334

$ws.Range("H93").Value = 21213.572
$ws.Range("I93").Value = 13299
$ws.Range("K93").Value = 13299
$ws.Range("M93").Value = -11427

$ws.Range("H95").Value = 51999.75
$ws.Range("J95").Value = 51999.75
$ws.Range("L95").Value = 51999.75
$ws.Range("N95").Value = -57491.75

$ws.Range("H103").Value = 39981
$ws.Range("I103").Value = 39981
$ws.Range("K103").Value = 39981
$ws.Range("M103").Value = -38809

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 111.92308
$ws.Range("I19").Value = 131
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 393
$ws.Range("L19").Value = 300
$ws.Range("M19").Value = -219
$ws.Range("N19").Value = -648

$ws.Range("H55").Value = 3975.652
$ws.Range("J55").Value = 4028.603
$ws.Range("L55").Value = 12085.809
$ws.Range("N55").Value = -12439.809

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null

$ws.Range("H103").Value = 386.66666
$ws.Range("I103").Value = 412.5
$ws.Range("K103").Value = 1237.5
$ws.Range("M103").Value = -358.5

$ws.Range("H114").Value = 1292.2727
$ws.Range("J114").Value = 1187.4445
$ws.Range("L114").Value = 3562.3335
$ws.Range("N114").Value = -10070.3335

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 14986.8
$ws.Range("I49").Value = 10000
$ws.Range("J49").Value = 26622.666
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 26622.666
$ws.Range("M49").Value = -9816
$ws.Range("N49").Value = -26990.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4995
$ws.Range("I7").Value = 4995
$ws.Range("K7").Value = 4995
$ws.Range("M7").Value = -4883

$ws.Range("H13").Value = 7500700
$ws.Range("I13").Value = 7500700
$ws.Range("K13").Value = 7500700
$ws.Range("M13").Value = -7500560

$ws.Range("H20").Value = 3002.5
$ws.Range("I20").Value = 3002.8572
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 3002.8572
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -2776.8572
$ws.Range("N20").Value = -3452

$ws.Range("H55").Value = 1285.4667
$ws.Range("I55").Value = 1067.2222
$ws.Range("J55").Value = 1612.8334
$ws.Range("K55").Value = 1067.2222
$ws.Range("L55").Value = 1612.8334
$ws.Range("M55").Value = -894.2221999999999
$ws.Range("N55").Value = -1958.8334

$ws.Range("H68").Value = 5270.75
$ws.Range("J68").Value = 6000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7498

$ws.Range("H71").Value = 5270.75
$ws.Range("J71").Value = 6000
$ws.Range("L71").Value = 30000
$ws.Range("N71").Value = -37488

$ws.Range("H93").Value = 66669290
$ws.Range("J93").Value = 2347.5
$ws.Range("L93").Value = 2347.5
$ws.Range("N93").Value = -4843.5

$ws.Range("H126").Value = 4995
$ws.Range("I126").Value = 4995
$ws.Range("K126").Value = 14985
$ws.Range("M126").Value = -12515

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 6750
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 3500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -3386
$ws.Range("N3").Value = -10228

$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10344

$ws.Range("H20").Value = 3505
$ws.Range("J20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("N20").Value = -7480

$ws.Range("H113").Value = 660.46155
$ws.Range("I113").Value = 540.6667
$ws.Range("K113").Value = 1622.0001
$ws.Range("M113").Value = 547.9999
